$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New localization rows for the "bonus round" feature.
# Columns: A=Key, B=Value, C=VoiceDuration
# Cell writes are ordered to match the original authoring sequence so the
# shared-strings table comes out in the same order.

$ws.Cells.Item(108, 1).Value = "bonus_round"
$ws.Cells.Item(108, 2).Value = "BONUS ROUND"
$ws.Cells.Item(108, 3).Value = 2

$ws.Cells.Item(109, 1).Value = "commutative"
$ws.Cells.Item(110, 1).Value = "associative"
$ws.Cells.Item(111, 1).Value = "distributive"

$ws.Cells.Item(109, 2).Value = "COMMUTATIVE"
$ws.Cells.Item(110, 2).Value = "ASSOCIATIVE"
$ws.Cells.Item(111, 2).Value = "DISTRIBUTIVE"

$ws.Cells.Item(109, 3).Value = 1
$ws.Cells.Item(110, 3).Value = 1
$ws.Cells.Item(111, 3).Value = 1

$ws.Cells.Item(112, 1).Value = "bonus_instruct"
$ws.Cells.Item(112, 2).Value = "Drag the correct numbers on the slots."
$ws.Cells.Item(112, 3).Value = 3

$ws.Cells.Item(113, 1).Value = "proceed"
$ws.Cells.Item(113, 2).Value = "PROCEED"

$ws.Cells.Item(115, 1).Value = "bonus_incorrect"
$ws.Cells.Item(115, 2).Value = "INCORRECT!"

$ws.Cells.Item(116, 1).Value = "bonus_score_format"
$ws.Cells.Item(116, 2).Value = "BONUS SCORE: +{0}"

$ws.Cells.Item(114, 1).Value = "bonus_time_expired"
$ws.Cells.Item(114, 2).Value = "TIME'S UP!"

$ws.Application.ActiveWindow.ScrollRow = 100
$ws.Range("A114").Select()
